$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue "D2" "285.79"
Set-TextValue "E2" "-10.43%"
Set-TextValue "D3" "39.89"
Set-TextValue "E3" "-3.75%"
Set-TextValue "D4" "5.042"
Set-TextValue "E4" "-3.77%"
Set-TextValue "D5" "0.07274"
Set-TextValue "E5" "-5.97%"
Set-TextValue "B6" "GateToken"
Set-TextValue "C6" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D6" "4.310"
Set-TextValue "E6" "-0.24%"
Set-TextValue "B7" "FTXToken"
Set-TextValue "C7" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D7" "1.510"
Set-TextValue "E7" "-11.06%"
Set-TextValue "B8" "MXToken"
Set-TextValue "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.9146"
Set-TextValue "E8" "-3.90%"
Set-TextValue "B9" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D9" "0.1203"
Set-TextValue "E9" "-4.59%"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1705"
Set-TextValue "E10" "-7.01%"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.08594"
Set-TextValue "E11" "-6.27%"
Set-TextValue "B12" "BitrueCoin"
Set-TextValue "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.04166"
Set-TextValue "E12" "-4.66%"
Set-TextValue "B13" "BitMartToken"
Set-TextValue "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.1050"
Set-TextValue "E13" "-0.02%"
Set-TextValue "B14" "BitForexToken"
Set-TextValue "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001276"
Set-TextValue "E14" "-0.92%"
Set-TextValue "B15" "TigerCash"
Set-TextValue "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.005932"
Set-TextValue "E15" "-0.93%"
Set-TextValue "B16" "LEO"
Set-TextValue "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.403"
Set-TextValue "E16" "1.97%"
Set-TextValue "D18" "0.3284"
Set-TextValue "E18" "-2.06%"
Set-TextValue "D19" "7.845"
Set-TextValue "E19" "1.87%"
Set-TextValue "D20" "0.1345"
Set-TextValue "E20" "-0.50%"
Set-TextValue "E22" "-4.42%"
Set-TextValue "D23" "0.001271"
Set-TextValue "E23" "0.64%"
Set-TextValue "D24" "0.003783"
Set-TextValue "E24" "-8.29%"
Set-TextValue "D25" "0.0001283"
Set-TextValue "E25" "1.09%"
Set-TextValue "D26" "0.0003735"
Set-TextValue "D38" "0.02283"
Set-TextValue "E38" "-10.70%"
Set-TextValue "D39" "0.04914"
Set-TextValue "E39" "-8.29%"
Set-TextValue "D40" "0.006797"
Set-TextValue "E40" "240.13%"
Set-TextValue "D41" "0.007685"
Set-TextValue "E41" "-0.96%"
Set-TextValue "D42" "0.1266"
Set-TextValue "E42" "-3.96%"
Set-TextValue "D43" "0.007407"
Set-TextValue "E43" "0.92%"
Set-TextValue "D44" "0.006930"
Set-TextValue "E44" "-8.50%"
Set-TextValue "D45" "0.3078"
Set-TextValue "E45" "-10.75%"
Set-TextValue "D46" "0.00006403"
Set-TextValue "E46" "-4.21%"
Set-TextValue "E47" "0.42%"
Set-TextValue "E48" "34.76%"
Set-TextValue "E49" "0.12%"
Set-TextValue "D50" "0.00002107"
Set-TextValue "E50" "0.42%"
Set-TextValue "D51" "0.0002007"
Set-TextValue "E51" "0.42%"
